# Applies the "Trade #29 closed" update to the live trading results workbook.
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1198.97
$summary.Range("B4").Value = -1.03
$summary.Range("B5").Value = -0.71
$summary.Range("B6").Value = 29
$summary.Range("B8").Value = 17
$summary.Range("B9").Value = 20.69

# --- Strategy Status sheet (MarketMaking strategy row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 98.97
$status.Range("D4").Value = 29
$status.Range("E4").Value = -1.03
$status.Range("F4").Value = -1.03
$status.Range("G4").Value = 20.69

# --- All Trades & MarketMaking sheets (Trade #29, row 30) ---
$tradeSheets = @("All Trades", "MarketMaking")
foreach ($sheetName in $tradeSheets) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $sheet.Range("G30").Value = 0.38
    $sheet.Range("H30").Value = "CLOSED"
    $sheet.Range("I30").Value = -5
    $sheet.Range("J30").Value = -0.02
    $sheet.Range("K30").Value = 98.97
    $sheet.Range("P30").Value = "early_exit"
    $sheet.Range("Q30").Value = 0.13
}

$wb.Save()
